# Adds an example plate layout for EchoProto.Templates.Loop_Assembly
# to the "Well lookup" sheet of the "Example DNA Stocks" workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Well lookup")

# Name (column D) for each well, keyed by row number (row 2 = well A1 ... row 37 = well B12)
$namesByRow = @{
    2  = "pOdd1"
    3  = "pOdd1"
    4  = "pOdd1"
    5  = "J23100"
    6  = "B0034"
    7  = "GFP"
    8  = "GFP"
    9  = "GFP"
    10 = "mCherry"
    11 = "mCherry"
    12 = "mCherry"
    13 = "B0015"
    14 = "B0015"
    15 = "B0015"
    16 = "J23119"
    17 = "J23101"
    18 = "J23102"
    19 = "J23103"
    20 = "J23104"
    21 = "J23105"
    22 = "J23106"
    23 = "J23107"
    24 = "J23108"
    25 = "J23109"
    26 = "J23110"
    27 = "J23111"
    28 = "J23112"
    29 = "J23113"
    30 = "J23114"
    31 = "J23115"
    32 = "J23116"
    33 = "J23117"
    34 = "J23118"
    35 = "B0030"
    36 = "B0031"
    37 = "B0032"
}

# Order in which the rows were originally populated
$rowOrder = @(2,3,4,5,6,7,8,9,10,11,12,13,14,15,17,18,16,19,20,21,22,23,24,25,26,27,28,29,30,31,32,33,34,35,36,37)

foreach ($row in $rowOrder) {
    $ws.Cells.Item($row, 4).Value = $namesByRow[$row]  # Column D - Name
    $ws.Cells.Item($row, 5).Value = 50                 # Column E - Volume (uL) - Initial
    $ws.Cells.Item($row, 9).Value = "AQ_BP"             # Column I - Calibration Type
}

# Restore the active selection as recorded after the edit
$ws.Range("D15").Select()
